# Generate Report for Archive
# -----------------------------------------------------------------------
# The nightly localization-status report moved the two in-flight docs
# ("9d52b4be-...md" / "ba7e476f-...md") from "Ready for handoff" to
# "In Translation" on both the per-locale sheets (Status column) and the
# Overview roll-up (one column per locale). Re-running the report also
# re-autosizes the now-narrower Status/locale columns so the shorter
# text doesn't leave a ragged gap.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: one status column per locale (zh-cn = E, de-de = F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Column shrinks to fit the shorter "In Translation" label.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale sheets: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
